# Actualizacion desde MV -datos-
# Appends new daily rows (04-08-2021 .. 02-09-2021) to Sheet1,
# matching the "Recompra deuda BCCh 2021 - Diaria" data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("04-08-2021", 899, 5115, 76)
    ,@("05-08-2021", 897, 5108, 76)
    ,@("06-08-2021", 896, 5101, 76)
    ,@("07-08-2021", 896, 5101, 76)
    ,@("08-08-2021", 896, 5101, 76)
    ,@("09-08-2021", 887, 5050, 75)
    ,@("10-08-2021", 884, 5031, 75)
    ,@("11-08-2021", 892, 5079, 76)
    ,@("12-08-2021", 900, 5125, 76)
    ,@("13-08-2021", 901, 5127, 76)
    ,@("14-08-2021", 901, 5127, 76)
    ,@("15-08-2021", 901, 5127, 76)
    ,@("16-08-2021", 899, 5117, 76)
    ,@("17-08-2021", 888, 5052, 75)
    ,@("18-08-2021", 881, 5017, 75)
    ,@("19-08-2021", 884, 5029, 75)
    ,@("20-08-2021", 881, 5015, 75)
    ,@("21-08-2021", 881, 5015, 75)
    ,@("22-08-2021", 881, 5015, 75)
    ,@("23-08-2021", 885, 5038, 75)
    ,@("24-08-2021", 889, 5058, 75)
    ,@("25-08-2021", 890, 5065, 75)
    ,@("26-08-2021", 889, 5060, 75)
    ,@("27-08-2021", 887, 5046, 75)
    ,@("28-08-2021", 887, 5046, 75)
    ,@("29-08-2021", 887, 5046, 75)
    ,@("30-08-2021", 887, 5051, 75)
    ,@("31-08-2021", 892, 5079, 76)
    ,@("01-09-2021", 898, 5111, 76)
    ,@("02-09-2021", 907, 5165, 77)
)

$startRow = 217
$rowCount = $data.Length
$helperCol = 26   # column Z, used as scratch space so Excel does not
                   # auto-convert the dd-mm-yyyy text into a real date

# Build helper formulas that evaluate to plain text dates, so the
# result can be pasted as values without triggering date recognition.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, $helperCol).Formula = '="' + $data[$i][0] + '"'
}

$helperRange = $ws.Range($ws.Cells.Item($startRow, $helperCol), $ws.Cells.Item($startRow + $rowCount - 1, $helperCol))
$helperRange.Copy()

$targetA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $rowCount - 1, 1))
$targetA.PasteSpecial(-4163)   # xlPasteValues

$helperRange.ClearContents()

# Fill in the numeric columns B, C and D.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
}
